$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Two new work-order rows (28 and 29, i.e. sheet rows 30 and 31) were added to
# the "Report" table, pushing the printed/used range from A1:AK29 to A1:AK31.
# ---------------------------------------------------------------------------

# 1) Clone the row formatting first, before any content is written, using the
#    most similar existing rows as templates:
#      - row 30 continues the "white" banding of row 28
#      - row 31 continues the "shaded" banding of row 29
$ws.Range('A28:AK28').Copy() | Out-Null
$ws.Range('A30:AK30').PasteSpecial(-4122) | Out-Null

$ws.Range('A29:AK29').Copy() | Out-Null
$ws.Range('A31:AK31').PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# 2) The "報修說明" (P) and "工作內容" (AC) cells on rows 29 and 30 switch to the
#    wrapping variant of their column style.
$ws.Range('P29').WrapText = $true
$ws.Range('AC29').WrapText = $true
$ws.Range('P30').WrapText = $true
$ws.Range('AC30').WrapText = $true

# 3) Row 30 (item 28) - 三重正義北店 / D028
$ws.Range('A30').Value2 = 28
$ws.Range('B30').Value2 = '服務'
$ws.Range('C30').Value2 = 2025070792
$ws.Range('F30').Value2 = 'D028'
$ws.Range('G30').Value2 = '三重正義北店'
$ws.Range('H30').Value2 = '新北市三重區'
$ws.Range('Q30').Value2 = 'THILF0D028'
$ws.Range('R30').Value2 = '新北一'
$ws.Range('S30').Value2 = '吳宗鴻'
$ws.Range('T30').Value2 = 1
$ws.Range('U30').Value2 = '已完工'
$ws.Range('V30').Value2 = '2025-07-04 15:38:39'
$ws.Range('W30').Value2 = '2025-07-04 15:00:00'
$ws.Range('X30').Value2 = '2025-07-04 15:38:00'
$ws.Range('Z30').Value2 = 0.6
$ws.Range('AB30').Value2 = '到場處理'
$ws.Range('AC30').Value2 = 'PMQ3+STAR'
$ws.Range('AD30').Value2 = 'O'
$ws.Range('AJ30').Value2 = 'O'
$ws.Range('AK30').Value2 = 'O'

# 4) Row 31 (item 29) - 北縣重武店 / 3362
$ws.Range('A31').Value2 = 29
$ws.Range('B31').Value2 = '服務'
$ws.Range('C31').Value2 = 2025070831
$ws.Range('F31').Value2 = 3362
$ws.Range('G31').Value2 = '北縣重武店'
$ws.Range('H31').Value2 = '新北市三重區'
$ws.Range('Q31').Value2 = 'THILF03362'
$ws.Range('R31').Value2 = '新北一'
$ws.Range('S31').Value2 = '吳宗鴻'
$ws.Range('T31').Value2 = 1
$ws.Range('U31').Value2 = '已完工'
$ws.Range('V31').Value2 = '2025-07-04 16:09:19'
$ws.Range('W31').Value2 = '2025-07-04 15:40:00'
$ws.Range('X31').Value2 = '2025-07-04 16:08:00'
$ws.Range('Z31').Value2 = 0.5
$ws.Range('AB31').Value2 = '到場處理'
$ws.Range('AC31').Value2 = 'PMQ3+STAR'
$ws.Range('AD31').Value2 = 'O'
$ws.Range('AJ31').Value2 = 'O'
$ws.Range('AK31').Value2 = 'O'

# 5) The print area grows by two rows to keep covering the whole table.
$ws.PageSetup.PrintArea = '$A$1:$AK$31'

# 6) Leave the cursor where the author left it after typing the last row.
$ws.Range('A31').Select()
